# Update MSME Lithuania Summary country-indicator figures with refreshed
# (more precise) percentages, per the updated SBS Eurostat source data.
#
# The affected cells hold numeric-looking values that are stored as TEXT
# in the workbook (shared strings), so we force a Text number format
# before writing the new value — otherwise Excel's automatic type
# inference would coerce them into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B13:D13").NumberFormat = "@"
$ws.Range("B13").Value = "37.86"
$ws.Range("C13").Value = "4.13"
$ws.Range("D13").Value = "41.99"

# Employment (% of total): Micro / SMEs / MSMEs
$ws.Range("B14:D14").NumberFormat = "@"
$ws.Range("B14").Value = "24.86"
$ws.Range("C14").Value = "51.54"
$ws.Range("D14").Value = "76.39"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B16:D16").NumberFormat = "@"
$ws.Range("B16").Value = "89.97"
$ws.Range("C16").Value = "9.81"
$ws.Range("D16").Value = "99.77"

# Value added to the economy (% of total): Micro / SMEs / MSMEs
$ws.Range("B20:D20").NumberFormat = "@"
$ws.Range("B20").Value = "13.48"
$ws.Range("C20").Value = "55.49"
$ws.Range("D20").Value = "68.98"
